$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the D5 cell value (u & v calculation fix)
$ws.Range("D5").Value = 500

# Update the selected cell/active selection to D5
$ws.Range("D5").Select()
